$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A labels for the two new rows first (this ordering matters for the
# shared-string table layout Excel produces on save).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# ---------------------------------------------------------------------------
# Row 2, column B: update the Cases-tab WebExcel query text (indentation of
# the WHERE clause + coalesce()-wrapped Age (years) expression).
# ---------------------------------------------------------------------------
$casesQuery = @"
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
        WHERE   d.tumor_grade IN ["Not Reported"] 
return ss.study_subject_id as ``Case ID``,
       p.program_acronym as ``Program Code``,
        p.program_id as Program_ID,
       s.study_acronym as ``Arm``,
       ss.disease_subtype as ``Diagnosis``,
       sf.grouped_recurrence_score AS ``Recurrence Score``,
       d.tumor_size_group AS ``tumor_size``,
       d.er_status AS ``ER Status``,
       d.pr_status AS ``PR Status``,
       coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '') AS ``Age (years)``,
demo.survival_time AS ``Survival (days)``
"@
$ws.Range("B2").Value = $casesQuery

# ---------------------------------------------------------------------------
# Row 3 (SamplesTab): WebExcel query
# ---------------------------------------------------------------------------
$samplesQuery = @"
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
 WHERE   d.tumor_grade IN ["Not Reported"] 
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS ``Sample ID``,
            ss.study_subject_id AS ``Case ID``,
            p.program_acronym AS ``Program Code``,
            s.study_acronym AS ``Arm``,
            ss.disease_subtype AS ``Diagnosis``,
            samp.tissue_type AS ``Tissue Type``,
            samp.composition AS ``Tissue Composition``,
            samp.sample_anatomic_site AS ``Sample Anatomic Site``,
            samp.method_of_sample_procurement AS ``Sample Procurement Method``
"@
$ws.Range("B3").Value = $samplesQuery

# ---------------------------------------------------------------------------
# Row 4 (FilesTab): WebExcel query
# ---------------------------------------------------------------------------
$filesQuery = @"
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
 WHERE   d.tumor_grade IN ["Not Reported"] 
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS ``File Name``,
    head(labels(samp)) AS ``Association``,
    f.file_description AS ``Description``,
    f.file_format AS ``File Format``,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS ``Program Code``,
    s.study_acronym AS ``Arm``,
    ss.study_subject_id AS ``Case ID``,
    samp.sample_id AS ``Sample ID``
    order by f.file_name
"@
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# Remaining columns for rows 3 & 4 (StatQuery / dbExcel / WebExcel file
# names) mirror row 2's values.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2

# ---------------------------------------------------------------------------
# Wrap text on the query / StatQuery columns (B & C) for all three data rows,
# matching the existing "Normal 2" wrap style used on B2/C2.
# ---------------------------------------------------------------------------
$ws.Range("B2:C4").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights to fit the taller wrapped text blocks.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# ---------------------------------------------------------------------------
# Column widths (best-fit-ish, matching the final layout).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.78
$ws.Columns.Item(2).ColumnWidth = 76.11
$ws.Columns.Item(3).ColumnWidth = 47.89
$ws.Columns.Item(4).ColumnWidth = 60.78
$ws.Columns.Item(5).ColumnWidth = 59.44

# ---------------------------------------------------------------------------
# Reset the view: select the whole sheet (clears the lingering topLeftCell /
# stale B2 selection left over from the previous save).
# ---------------------------------------------------------------------------
$ws.Cells.Select() | Out-Null
